# Append a new scraped case ("Power Automate for Desktop ...") at the top of the
# existing "案件" list (row 5) and push the previously-captured rows down by one,
# refreshing the capture timestamp in column A and re-numbering the URL
# hyperlinks that live in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-26 01:22:50"

# --- Column B width: 47 -> 50 -------------------------------------------------
# ColumnWidth uses Excel's "characters" unit, which differs from the raw width
# stored in the sheet XML by Excel's standard padding (~0.91666 chars at the
# default font). Using 49.1667 here yields a stored width of exactly 50.
$ws.Columns.Item(2).ColumnWidth = 49.1667

# --- Capture the existing data rows (old rows 5-13) before overwriting -------
$captured = @()
for ($r = 5; $r -le 13; $r++) {
    $row = @{
        B = $ws.Cells.Item($r, 2).Text
        C = $ws.Cells.Item($r, 3).Text
        D = $ws.Cells.Item($r, 4).Text
        E = $ws.Cells.Item($r, 5).Text
        F = $ws.Cells.Item($r, 6).Text
        G = $ws.Cells.Item($r, 7).Text
        H = $ws.Cells.Item($r, 8).Text
    }
    $captured += $row
}

# --- Refresh the timestamp on the rows that are kept unchanged (rows 2-4) ----
for ($r = 2; $r -le 4; $r++) {
    $ws.Range("A$r").Value = $newTimestamp
}

# --- Remove all existing hyperlinks; they get rebuilt after the cell values
#     have been written into their final positions --------------------------
$ws.Hyperlinks.Delete()

# --- Write the brand-new row 5 ------------------------------------------------
$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "【Power Automate for Desktop】販売管理システムへExcelから自動入力"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5407216"
$ws.Range("G5").Value = 48
$ws.Range("H5").Value = "◇管理"

# --- Write back the captured rows, shifted down by one (rows 6-14) ----------
for ($i = 0; $i -le 8; $i++) {
    $r = 6 + $i
    $data = $captured[$i]
    $ws.Range("A$r").Value = $newTimestamp
    $ws.Range("B$r").Value = $data.B
    $ws.Range("C$r").Value = $data.C
    $ws.Range("D$r").Value = $data.D
    $ws.Range("E$r").Value = $data.E
    $ws.Range("F$r").Value = $data.F
    if ($data.G -ne "") {
        $ws.Range("G$r").Value = [double]$data.G
    }
    if ($data.H -ne "") {
        $ws.Range("H$r").Value = $data.H
    }
}

# --- Rebuild the column F hyperlinks (address = the cell's own URL text),
#     in row order, so relationship ids come out sequential again -----------
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $cell.Text)
    $cell.Style = "Hyperlink"
}
